$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.150.34'
$ws.Range('E2').Value = '  +2.98%  '

$ws.Range('D3').Value = '1.580.09'
$ws.Range('E3').Value = '  +1.80%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.78%  '

$ws.Range('E7').Value = '  -0.22%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.34'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.95%  '

$ws.Range('E9').Value = '  +2.37%  '

$ws.Range('E10').Value = '  +1.70%  '

$ws.Range('E11').Value = '  +1.60%  '

$ws.Range('D12').Value = '1.803.75'
$ws.Range('E12').Value = '  +1.67%  '

$ws.Range('D13').Value = '1.573.50'
$ws.Range('E13').Value = '  +1.46%  '

$ws.Range('D14').Value = '29.161.57'
$ws.Range('E14').Value = '  +3.10%  '

$ws.Range('E15').Value = '  +2.81%  '

$ws.Range('E16').Value = '  +2.72%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.83%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '236.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.71%  '

$ws.Range('E19').Value = '  +1.58%  '

$ws.Range('D20').Value = '0.0₃0691'
$ws.Range('E20').Value = '  +2.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '

$ws.Range('E22').Value = '  +1.81%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.21%  '

$ws.Range('E24').Value = '  +4.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.42%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.53%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.108'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.14%  '

$ws.Range('E28').Value = '  +1.71%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('E30').Value = '  +0.31%  '

$ws.Range('E31').Value = '  +0.28%  '

$ws.Range('E32').Value = '  +1.49%  '

$ws.Range('D33').Value = '1.422.14'

$ws.Range('E34').Value = '  +1.39%  '

$ws.Range('E35').Value = '  -1.88%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.04%  '

$ws.Range('E38').Value = '  -1.78%  '

$ws.Range('E39').Value = '  +1.12%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.529'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.61%  '

$ws.Range('E41').Value = '  +2.11%  '

$ws.Range('E42').Value = '  -0.16%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '52.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +24.85%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.789'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.43%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0469'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.45%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.55%  '

$ws.Range('E47').Value = '  -0.76%  '

$ws.Range('D48').Value = '1.716.90'
$ws.Range('E48').Value = '  +1.65%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.847'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.54%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.05%  '

$ws.Range('E51').Value = '  +1.64%  '
